$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '27.989.97'
$ws.Range("E2").Value = '  +7.17%  '

$ws.Range("D3").Value = '1.741.07'
$ws.Range("E3").Value = '  +5.35%  '

$ws.Range("E4").Value = '  -0.13%  '

Set-TextValue "D5" '227.99'

Set-TextValue "D6" '0.5430'
$ws.Range("E6").Value = '  +3.46%  '

$ws.Range("E7").Value = '  -0.15%  '

Set-TextValue "D8" '0.2761'
$ws.Range("E8").Value = '  +3.73%  '

Set-TextValue "D9" '0.06729'
$ws.Range("E9").Value = '  +6.03%  '

Set-TextValue "D10" '21.66'
$ws.Range("E10").Value = '  +4.98%  '

Set-TextValue "D11" '0.07790'
$ws.Range("E11").Value = '  +0.97%  '

Set-TextValue "D12" '4.699'
$ws.Range("E12").Value = '  +2.13%  '

$ws.Range("D13").Value = '1.743.52'
$ws.Range("E13").Value = '  +4.49%  '

$ws.Range("D14").Value = '1.980.94'
$ws.Range("E14").Value = '  +5.43%  '

Set-TextValue "D15" '0.5967'
$ws.Range("E15").Value = '  +6.21%  '

$ws.Range("D16").Value = '0.0₅8363'
$ws.Range("E16").Value = '  +2.10%  '

Set-TextValue "D17" '68.83'
$ws.Range("E17").Value = '  +5.30%  '

$ws.Range("D18").Value = '27.976.14'
$ws.Range("E18").Value = '  +7.11%  '

Set-TextValue "D19" '224.11'
$ws.Range("E19").Value = '  +17.24%  '

Set-TextValue "D20" '4.836'
$ws.Range("E20").Value = '  +2.84%  '

$ws.Range("E21").Value = '  -0.22%  '

$ws.Range("E22").Value = '  +5.13%  '

Set-TextValue "D23" '6.230'
$ws.Range("E23").Value = '  +4.08%  '

$ws.Range("E24").Value = '  -0.16%  '

Set-TextValue "D25" '146.15'
$ws.Range("E25").Value = '  +0.61%  '

$ws.Range("E26").Value = '  +3.46%  '

Set-TextValue "D27" '17.26'
$ws.Range("E27").Value = '  +8.37%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D28" '1.665'
$ws.Range("E28").Value = '  +10.40%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D29" '7.447'
$ws.Range("E29").Value = '  +2.52%  '

Set-TextValue "D30" '0.05642'
$ws.Range("E30").Value = '  +0.21%  '

Set-TextValue "D31" '1.315'
$ws.Range("E31").Value = '  +3.34%  '

Set-TextValue "D32" '3.700'
$ws.Range("E32").Value = '  +5.90%  '

Set-TextValue "D33" '3.515'
$ws.Range("E33").Value = '  +4.27%  '

Set-TextValue "D34" '1.673'
$ws.Range("E34").Value = '  +5.99%  '

$ws.Range("E35").Value = '  +3.15%  '

$ws.Range("E36").Value = '  +2.42%  '

Set-TextValue "D37" '2.454'
$ws.Range("E37").Value = '  +1.79%  '

$ws.Range("E38").Value = '  +3.51%  '

Set-TextValue "D39" '0.01663'
$ws.Range("E39").Value = '  +4.29%  '

Set-TextValue "D40" '5.943'
$ws.Range("E40").Value = '  -0.78%  '

Set-TextValue "D41" '0.8500'
$ws.Range("E41").Value = '  +1.37%  '

$ws.Range("D42").Value = '1.048.23'
$ws.Range("E42").Value = '  +3.33%  '

$ws.Range("E43").Value = '  -0.11%  '

Set-TextValue "D44" '102.11'
$ws.Range("E44").Value = '  +0.35%  '

$ws.Range("D45").Value = '1.886.50'
$ws.Range("E45").Value = '  +5.28%  '

$ws.Range("D46").Value = '0.0₈117'
$ws.Range("E46").Value = '  +12.21%  '

Set-TextValue "D47" '59.89'
$ws.Range("E47").Value = '  +2.73%  '

Set-TextValue "D48" '8.279'
$ws.Range("E48").Value = '  +3.70%  '

Set-TextValue "D49" '0.4435'
$ws.Range("E49").Value = '  +2.07%  '

Set-TextValue "D50" '1.001'
$ws.Range("E50").Value = '  -0.10%  '

Set-TextValue "D51" '0.05319'
$ws.Range("E51").Value = '  -0.18%  '
